$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Beta) values that changed
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 9.178312775821793
$ws.Range("G2").Value = 8.766803171361049
$ws.Range("H2").Value = 9.594480817265291
$ws.Range("I2").Value = 0.01538991438918501
$ws.Range("J2").Value = 0.01424563830712567
$ws.Range("K2").Value = 0.01672905441475471
$ws.Range("L2").Value = 0.005717097912475402
$ws.Range("M2").Value = 0.005396440862554793
$ws.Range("N2").Value = 0.006077740589567102

# Update row 3 (Gamma) values
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.277040445070052
$ws.Range("G3").Value = 0.2761937718556184
$ws.Range("H3").Value = 0.2779249621772381
$ws.Range("I3").Value = 0.2491529967964021
$ws.Range("J3").Value = 0.2483557484447161
$ws.Range("K3").Value = 0.2499856937723305
$ws.Range("L3").Value = 0.2744480533082039
$ws.Range("M3").Value = 0.2736050499399066
$ws.Range("N3").Value = 0.2753294172235052

# Add new row 4 (Beta + Gamma) - copy row 3 formats down first so A4 picks up
# the same style (s="1") as A2/A3, then overwrite with the new values.
$ws.Range("A3:N3").Copy()
$ws.Range("A4:N4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 9.455353220891844
$ws.Range("G4").Value = 9.042996943216668
$ws.Range("H4").Value = 9.872405779442531
$ws.Range("I4").Value = 0.2645429111855871
$ws.Range("J4").Value = 0.2626013867518417
$ws.Range("K4").Value = 0.2667147481870852
$ws.Range("L4").Value = 0.2801651512206793
$ws.Range("M4").Value = 0.2790014908024613
$ws.Range("N4").Value = 0.2814071578130723
